$wb = $excel.ActiveWorkbook

# --- Queries sheet: the "Aux State" threshold query text was rewritten ---
$wsQueries = $wb.Worksheets.Item("Queries")
$newAuxStateQuery = "SELECT [OldStatus] AS [Aux State],[dbo].[SECONDSTOhhmmss](SUM(LEFT([OldStatusTimeSpend],2) * 3600 + SUBSTRING([OldStatusTimeSpend], 4,2) * 60 + SUBSTRING([OldStatusTimeSpend], 7,2))) AS [Total Time Spent],`nFORMAT(MAX([InsertedDateTime]),'dd/MM/yyyy HH:mm:ss') as [Last Threshold Date Time],COUNT(1) AS [Threshold Count]`nFROM [dbo].[TDM_Alerts] WHERE [InsertedDateTime]>='ReportBeforeDate' AND [InsertedDateTime]<='ReportAfterDate' `nAND [AgentName] like 'AgentNameCapturedFromUI'`nGROUP BY [AgentName],[OldStatus], [TeamName] ORDER BY [Threshold Count] DESC,[Last Threshold Date Time] DESC;"
$wsQueries.Range("G2").Value = $newAuxStateQuery

# --- ShowDateRange sheet: report end date moved from 30-05-2020 to 11-11-2020 ---
$wsShowDateRange = $wb.Worksheets.Item("ShowDateRange")
# Use a leading apostrophe so the cell keeps being stored as (quote-prefixed) text,
# matching the original "Starts as text" formatting of these cells.
$wsShowDateRange.Range("E2:E7").Value = "'11-11-2020 00:00:00"

# --- Switch the active/selected sheet from Queries to AdvanceSearch ---
$wsAdvanceSearch = $wb.Worksheets.Item("AdvanceSearch")
$wsAdvanceSearch.Activate()
